$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# Rows 13,14,17,18,20,21 currently use cell style index 1 (alternate shading);
# in the target state every data row in the 9-21 block uses style index 2.
# Copy the formatting from a row that already has style 2 (row 9) onto those rows.
$ws.Range("A9:H9").Copy()
$ws.Range("A13:H13").PasteSpecial(-4122)
$ws.Range("A14:H14").PasteSpecial(-4122)
$ws.Range("A17:H17").PasteSpecial(-4122)
$ws.Range("A18:H18").PasteSpecial(-4122)
$ws.Range("A20:H20").PasteSpecial(-4122)
$ws.Range("A21:H21").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the "bus insert" panel rows (9-21) with labels/values/textbox names.
# The assignment order below matches the order the unique strings were
# originally entered in, so the shared-string table comes out in the same
# sequence (new strings are appended the first time each distinct value is
# written to a cell).
$ws.Range("A9").Value = "label7"
$ws.Range("B9").Value = "case"
$ws.Range("C9").Value = "combobox3"
$ws.Range("E9").Value = "trans/bus/insert"
$ws.Range("A10").Value = "label8"
$ws.Range("B10").Value = "Número da barra"
$ws.Range("C10").Value = "textbox1"
$ws.Range("A11").Value = "label9"
$ws.Range("B11").Value = "Número sequencial"
$ws.Range("C11").Value = "textbox8"
$ws.Range("A12").Value = "label10"
$ws.Range("B12").Value = "magnitude de tensão"
$ws.Range("C12").Value = "textbox2"
$ws.Range("A13").Value = "label11"
$ws.Range("B13").Value = "ângulo de fase"
$ws.Range("C13").Value = "textbox12"
$ws.Range("B14").Value = "base de tensão"
$ws.Range("A14").Value = "label12"
$ws.Range("C14").Value = "textbox3"
$ws.Range("B15").Value = "tensão especificada"
$ws.Range("A15").Value = "label16"
$ws.Range("C15").Value = "textbox11"
$ws.Range("B16").Value = "lim max geração"
$ws.Range("A16").Value = "label19"
$ws.Range("C16").Value = "textbox4"
$ws.Range("A18").Value = "label20"
$ws.Range("B17").Value = "lim min. geração"
$ws.Range("A19").Value = "label22"
$ws.Range("C19").Value = "textbox9"
$ws.Range("A17").Value = "label21"
$ws.Range("C17").Value = "textbox7"
$ws.Range("B18").Value = "lim max tensão"
$ws.Range("C18").Value = "textbox10"
$ws.Range("B19").Value = "lim min. tensão"
$ws.Range("A20").Value = "label23"
$ws.Range("B20").Value = "nome da barra"
$ws.Range("C20").Value = "textbox13"
$ws.Range("A21").Value = "label24"
$ws.Range("B21").Value = "area"
$ws.Range("C21").Value = "combobox4"

# Column D: plain numeric literal (6) for every row in the block.
$ws.Range("D9").Value = 6
$ws.Range("D10").Value = 6
$ws.Range("D11").Value = 6
$ws.Range("D12").Value = 6
$ws.Range("D13").Value = 6
$ws.Range("D14").Value = 6
$ws.Range("D15").Value = 6
$ws.Range("D16").Value = 6
$ws.Range("D17").Value = 6
$ws.Range("D18").Value = 6
$ws.Range("D19").Value = 6
$ws.Range("D20").Value = 6
$ws.Range("D21").Value = 6

# Column E: repeats the same shared string for the rest of the block.
$ws.Range("E10").Value = "trans/bus/insert"
$ws.Range("E11").Value = "trans/bus/insert"
$ws.Range("E12").Value = "trans/bus/insert"
$ws.Range("E13").Value = "trans/bus/insert"
$ws.Range("E14").Value = "trans/bus/insert"
$ws.Range("E15").Value = "trans/bus/insert"
$ws.Range("E16").Value = "trans/bus/insert"
$ws.Range("E17").Value = "trans/bus/insert"
$ws.Range("E18").Value = "trans/bus/insert"
$ws.Range("E19").Value = "trans/bus/insert"
$ws.Range("E20").Value = "trans/bus/insert"
$ws.Range("E21").Value = "trans/bus/insert"

# Update the active selection shown when the sheet is re-opened.
$ws.Range("D7").Select()
